# Update scripts with new TPM values.
# The "Sending cluster" (column A) changes from ECs to MuSCs for all rows.
# The "Target cluster" (column D) text stays the same for each row
# (ECs / FAPs / MuSCs respectively); only the computed numeric metrics
# in columns G, H, M, N, O, P, Q, R, S, T change with the new TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster (column A) changes from ECs to MuSCs for all data rows
$ws.Range("A2").Value = "MuSCs"
$ws.Range("A3").Value = "MuSCs"
$ws.Range("A4").Value = "MuSCs"

# Row 2 (Target cluster: ECs)
$ws.Range("G2").Value = 0.03636833333333334
$ws.Range("H2").Value = 0.109105
$ws.Range("M2").Value = 0.6882290000000001
$ws.Range("N2").Value = 2.064687
$ws.Range("O2").Value = 0.2885525922449623
$ws.Range("P2").Value = 0.2885525922449623
$ws.Range("Q2").Value = 0.02502974168166667
$ws.Range("R2").Value = 0.225267675135
$ws.Range("S2").Value = 0.2885525922449623
$ws.Range("T2").Value = 0.2885525922449623

# Row 3 (Target cluster: FAPs)
$ws.Range("G3").Value = 0.03636833333333334
$ws.Range("H3").Value = 0.109105
$ws.Range("O3").Value = 0.3570594926322683
$ws.Range("P3").Value = 0.3570594926322683
$ws.Range("Q3").Value = 0.03097219399777778
$ws.Range("R3").Value = 0.27874974598
$ws.Range("S3").Value = 0.3570594926322683
$ws.Range("T3").Value = 0.3570594926322683

# Row 4 (Target cluster: MuSCs)
$ws.Range("G4").Value = 0.03636833333333334
$ws.Range("H4").Value = 0.109105
$ws.Range("M4").Value = 0.8452533333333333
$ws.Range("N4").Value = 2.53576
$ws.Range("O4").Value = 0.3543879151227694
$ws.Range("P4").Value = 0.3543879151227694
$ws.Range("Q4").Value = 0.03074045497777778
$ws.Range("R4").Value = 0.2766640948
$ws.Range("S4").Value = 0.3543879151227694
$ws.Range("T4").Value = 0.3543879151227694
